$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-08 Sunday" "2024-12-09 Monday"

Replace-Text "75×74=5550" "43×70=3010"
Replace-Text "67×90=6030" "37×34=1258"
Replace-Text "71×55=3905" "53×19=1007"
Replace-Text "56×36=2016" "81×79=6399"
Replace-Text "26×26=676" "41×81=3321"
Replace-Text "65×31=2015" "33×20=660"
Replace-Text "30×54=1620" "70×41=2870"
Replace-Text "19×46=874" "56×83=4648"
Replace-Text "63×74=4662" "97×23=2231"
Replace-Text "83×79=6557" "20×22=440"
Replace-Text "29×94=2726" "23×93=2139"
Replace-Text "21×83=1743" "90×74=6660"
Replace-Text "20×83=1660" "16×95=1520"
Replace-Text "68×44=2992" "91×87=7917"
Replace-Text "91×28=2548" "77×32=2464"
Replace-Text "70×42=2940" "18×66=1188"
Replace-Text "42×29=1218" "27×21=567"
Replace-Text "38×79=3002" "92×68=6256"
Replace-Text "37×66=2442" "24×40=960"
Replace-Text "80×15=1200" "17×85=1445"
Replace-Text "71×12=852" "16×56=896"
Replace-Text "48×13=624" "99×82=8118"
Replace-Text "87×15=1305" "85×34=2890"
Replace-Text "51×63=3213" "24×37=888"
Replace-Text "95×38=3610" "44×16=704"
